$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell is written as text: force Text number format before the
# assignment so numeric-looking strings (e.g. "141.70", "1.00") are not
# auto-converted to numbers, then restore the "Normal" style so no extra
# formatting is left behind on the cell.
function Set-TextValue([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '63.596.16'
Set-TextValue 'E2' '  -1.74%  '
Set-TextValue 'D3' '3.042.08'
Set-TextValue 'E3' '  -2.00%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '557.26'
Set-TextValue 'E5' '  -0.48%  '
Set-TextValue 'D6' '141.70'
Set-TextValue 'E6' '  -1.78%  '
Set-TextValue 'E7' '  +0.03%  '
Set-TextValue 'D8' '3.039.10'
Set-TextValue 'E8' '  -1.87%  '
Set-TextValue 'D9' '0.517'
Set-TextValue 'E9' '  +3.19%  '
Set-TextValue 'B10' 'Toncoin'
Set-TextValue 'C10' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D10' '6.33'
Set-TextValue 'E10' '  -11.54%  '
Set-TextValue 'B11' 'Dogecoin'
Set-TextValue 'C11' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D11' '0.152'
Set-TextValue 'E11' '  +0.05%  '
Set-TextValue 'E12' '  +5.46%  '
Set-TextValue 'E13' '  +0.29%  '
Set-TextValue 'D14' '35.53'
Set-TextValue 'E14' '  +0.22%  '
Set-TextValue 'D15' '3.536.53'
Set-TextValue 'E15' '  -1.96%  '
Set-TextValue 'D16' '63.636.82'
Set-TextValue 'E16' '  -1.64%  '
Set-TextValue 'D17' '3.040.82'
Set-TextValue 'E17' '  -2.11%  '
Set-TextValue 'E18' '  +0.28%  '
Set-TextValue 'E19' '  -0.52%  '
Set-TextValue 'D20' '473.82'
Set-TextValue 'E20' '  -2.22%  '
Set-TextValue 'D21' '14.02'
Set-TextValue 'E21' '  +1.23%  '
Set-TextValue 'D22' '14.56'
Set-TextValue 'E22' '  +9.43%  '
Set-TextValue 'E23' '  +0.61%  '
Set-TextValue 'D24' '7.49'
Set-TextValue 'E24' '  -2.24%  '
Set-TextValue 'D25' '82.59'
Set-TextValue 'E25' '  +1.89%  '
Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  +0.00%  '
Set-TextValue 'E27' '  -0.79%  '
Set-TextValue 'D28' '8.08'
Set-TextValue 'E28' '  -0.44%  '
Set-TextValue 'E29' '  -2.41%  '
Set-TextValue 'D30' '0.999'
Set-TextValue 'E30' '  +0.03%  '
Set-TextValue 'D31' '26.14'
Set-TextValue 'E31' '  -0.05%  '
Set-TextValue 'E32' '  -1.54%  '
Set-TextValue 'E33' '  -1.48%  '
Set-TextValue 'E34' '  -0.26%  '
Set-TextValue 'E35' '  +0.16%  '
Set-TextValue 'D36' '54.67'
Set-TextValue 'E36' '  -1.18%  '
Set-TextValue 'E37' '  -0.69%  '
Set-TextValue 'D38' '440.72'
Set-TextValue 'E38' '  -5.57%  '
Set-TextValue 'D39' '0.0810'
Set-TextValue 'E39' '  -2.35%  '
Set-TextValue 'B40' 'dogwifhat'
Set-TextValue 'C40' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D40' '2.79'
Set-TextValue 'E40' '  +2.33%  '
Set-TextValue 'B41' 'Maker'
Set-TextValue 'C41' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D41' '3.007.34'
Set-TextValue 'E41' '  -0.41%  '
Set-TextValue 'E42' '  +0.29%  '
Set-TextValue 'D43' '8.24'
Set-TextValue 'E43' '  -0.61%  '
Set-TextValue 'D44' '0.268'
Set-TextValue 'E44' '  +2.24%  '
Set-TextValue 'D45' '27.65'
Set-TextValue 'E45' '  -3.33%  '
Set-TextValue 'D46' '2.23'
Set-TextValue 'E46' '  +6.32%  '
Set-TextValue 'E47' '  -0.03%  '
Set-TextValue 'E48' '  +0.48%  '
Set-TextValue 'D49' '118.29'
Set-TextValue 'E49' '  -0.50%  '
Set-TextValue 'E50' '  -1.04%  '
Set-TextValue 'E51' '  +0.01%  '
